$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3942060214836545
$ws.Range("C2").Value = 0.0324612055215141
$ws.Range("E2").Value = 0.3736801726878198
$ws.Range("F2").Value = 2.484156628171874
$ws.Range("G2").Value = 0.9701951879715125
$ws.Range("H2").Value = 0.9989775711007809
$ws.Range("J2").Value = 0.07164481598474115
$ws.Range("K2").Value = 0.3497029370666382
$ws.Range("M2").Value = 0.4015363223403909
$ws.Range("N2").Value = 1.962409002470153
$ws.Range("B3").Value = 0.3610762776451395
$ws.Range("C3").Value = 0.0283801183380632
$ws.Range("E3").Value = 0.3667164599643584
$ws.Range("F3").Value = 2.46711159199144
$ws.Range("G3").Value = 0.9702683033993793
$ws.Range("H3").Value = 1.003658630097249
$ws.Range("J3").Value = 0.07228467223535873
$ws.Range("K3").Value = 0.315127941958508
$ws.Range("M3").Value = 0.3835852324529441
$ws.Range("N3").Value = 1.98247249876113
$ws.Range("B4").Value = 0.340887356978044
$ws.Range("C4").Value = 0.02586522580301676
$ws.Range("E4").Value = 0.3626099348523226
$ws.Range("F4").Value = 2.457950600409802
$ws.Range("G4").Value = 0.9708998785524443
$ws.Range("H4").Value = 1.006963463108775
$ws.Range("J4").Value = 0.07269834371115769
$ws.Range("K4").Value = 0.2940010523925451
$ws.Range("M4").Value = 0.3727590154562108
$ws.Range("N4").Value = 1.995439591029999
$ws.Range("B5").Value = 0.3326988623188925
$ws.Range("C5").Value = 0.02483810421168187
$ws.Range("E5").Value = 0.3609791079356626
$ws.Range("F5").Value = 2.454545457449626
$ws.Range("G5").Value = 0.9713045748818701
$ws.Range("H5").Value = 1.00841850649141
$ws.Range("J5").Value = 0.07287214919222018
$ws.Range("K5").Value = 0.2854175723223022
$ws.Range("M5").Value = 0.3683965988942148
$ws.Range("N5").Value = 2.000886655846905
$ws.Range("B6").Value = 0.3313415099620158
$ws.Range("C6").Value = 0.02466741369308068
$ws.Range("E6").Value = 0.3607108861072987
$ws.Range("F6").Value = 2.453999847271845
$ws.Range("G6").Value = 0.9713806665561862
$ws.Range("H6").Value = 1.008666656484422
$ws.Range("J6").Value = 0.07290132544111749
$ws.Range("K6").Value = 0.2839938614572901
$ws.Range("M6").Value = 0.3676752072174594
$ws.Range("N6").Value = 2.001800971383375
$ws.Range("B7").Value = 0.3407767672817101
$ws.Range("C7").Value = 0.02585138292153033
$ws.Range("E7").Value = 0.3625877683180505
$ws.Range("F7").Value = 2.457903349326898
$ws.Range("G7").Value = 0.9709047401876205
$ws.Range("H7").Value = 1.006982647833397
$ws.Range("J7").Value = 0.07270066653051499
$ws.Range("K7").Value = 0.2938851874803987
$ws.Range("M7").Value = 0.37269998239465
$ws.Range("N7").Value = 1.995512392711408
$ws.Range("B8").Value = 0.382751272776062
$ws.Range("C8").Value = 0.03105593348503533
$ws.Range("E8").Value = 0.3712439970607804
$ws.Range("F8").Value = 2.478008660675428
$ws.Range("G8").Value = 0.9700985080008735
$ws.Range("H8").Value = 1.000502223823872
$ws.Range("J8").Value = 0.07186112285842805
$ws.Range("K8").Value = 0.3377603409381322
$ws.Range("M8").Value = 0.395306204289966
$ws.Range("N8").Value = 1.969192351206075
$ws.Range("B9").Value = 0.4662713171084647
$ws.Range("C9").Value = 0.04119029715448619
$ws.Range("E9").Value = 0.389560063378525
$ws.Range("F9").Value = 2.527795256360761
$ws.Range("G9").Value = 0.9731829793276745
$ws.Range("H9").Value = 0.9912112799397619
$ws.Range("J9").Value = 0.07037971122753506
$ws.Range("K9").Value = 0.4246079302608337
$ws.Range("M9").Value = 0.4411884769654932
$ws.Range("N9").Value = 1.922723503694066
$ws.Range("B10").Value = 0.5283701808154717
$ws.Range("C10").Value = 0.04859350886988523
$ws.Range("E10").Value = 0.4038344510296668
$ws.Range("F10").Value = 2.570707848299108
$ws.Range("G10").Value = 0.9783097179580977
$ws.Range("H10").Value = 0.9864694619441394
$ws.Range("J10").Value = 0.06939171421968515
$ws.Range("K10").Value = 0.4889103176485889
$ws.Range("M10").Value = 0.4758450642645613
$ws.Range("N10").Value = 1.89172204259901
$ws.Range("B11").Value = 0.5567809835802962
$ws.Range("C11").Value = 0.05195250156413067
$ws.Range("E11").Value = 0.4105058818908063
$ws.Range("F11").Value = 2.591609895431532
$ws.Range("G11").Value = 0.9812668191596003
$ws.Range("H11").Value = 0.9847651696293838
$ws.Range("J11").Value = 0.06896401707358679
$ws.Range("K11").Value = 0.518271697426087
$ws.Range("M11").Value = 0.4918173294867145
$ws.Range("N11").Value = 1.878300957123088
$ws.Range("B12").Value = 0.5675625564562097
$ws.Range("C12").Value = 0.05322321309591871
$ws.Range("E12").Value = 0.4130577355225427
$ws.Range("F12").Value = 2.599723745692472
$ws.Range("G12").Value = 0.9824767242135835
$ws.Range("H12").Value = 0.9841849268745761
$ws.Range("J12").Value = 0.06880518543839376
$ws.Range("K12").Value = 0.5294058177885859
$ws.Range("M12").Value = 0.4978952999078388
$ws.Range("N12").Value = 1.873316859361708
$ws.Range("B13").Value = 0.565239531836113
$ws.Range("C13").Value = 0.05294959903012852
$ws.Range("E13").Value = 0.4125070137914264
$ws.Range("F13").Value = 2.597967443623219
$ws.Range("G13").Value = 0.9822121375999586
$ws.Range("H13").Value = 0.9843069951774339
$ws.Range("J13").Value = 0.06883925354246001
$ws.Range("K13").Value = 0.5270071949895794
$ws.Range("M13").Value = 0.4965849847473649
$ws.Range("N13").Value = 1.874385905033925
$ws.Range("B14").Value = 0.5576675300978025
$ws.Range("C14").Value = 0.05205706917311659
$ws.Range("E14").Value = 0.410715313157219
$ws.Range("F14").Value = 2.592273443817859
$ws.Range("G14").Value = 0.9813645512910512
$ws.Range("H14").Value = 0.984716127138384
$ws.Range("J14").Value = 0.06895088720504061
$ws.Range("K14").Value = 0.5191873966471974
$ws.Range("M14").Value = 0.4923167748164445
$ws.Range("N14").Value = 1.877888943890568
$ws.Range("B15").Value = 0.553032451071374
$ws.Range("C15").Value = 0.05151020381632065
$ws.Range("E15").Value = 0.409621167450382
$ws.Range("F15").Value = 2.58881158479204
$ws.Range("G15").Value = 0.9808571233314183
$ws.Range("H15").Value = 0.9849752156657843
$ws.Range("J15").Value = 0.06901967342665438
$ws.Range("K15").Value = 0.5143995704081021
$ws.Range("M15").Value = 0.4897062269026264
$ws.Range("N15").Value = 1.88004744649961
$ws.Range("B16").Value = 0.5265167020005208
$ws.Range("C16").Value = 0.04837381392738394
$ws.Range("E16").Value = 0.4034020328051682
$ws.Range("F16").Value = 2.569369648340682
$ws.Range("G16").Value = 0.9781290609876265
$ws.Range("H16").Value = 0.986589954822179
$ws.Range("J16").Value = 0.06942010291038248
$ws.Range("K16").Value = 0.486993676760278
$ws.Range("M16").Value = 0.4748053905471892
$ws.Range("N16").Value = 1.892612871737068
$ws.Range("B17").Value = 0.5102914003186925
$ws.Range("C17").Value = 0.04644749121386837
$ws.Range("E17").Value = 0.3996323270707691
$ws.Range("F17").Value = 2.557796420205833
$ws.Range("G17").Value = 0.9766157196258263
$ws.Range("H17").Value = 0.9876965320076607
$ws.Range("J17").Value = 0.06967132399611886
$ws.Range("K17").Value = 0.470209054090958
$ws.Range("M17").Value = 0.4657170943761813
$ws.Range("N17").Value = 1.900496047259743
$ws.Range("B18").Value = 0.5009742717145969
$ws.Range("C18").Value = 0.04533869662007817
$ws.Range("E18").Value = 0.3974808441888911
$ws.Range("F18").Value = 2.551269773925696
$ws.Range("G18").Value = 0.9758040911919039
$ws.Range("H18").Value = 0.9883756225171538
$ws.Range("J18").Value = 0.06981786666279044
$ws.Range("K18").Value = 0.4605653394444005
$ws.Range("M18").Value = 0.4605092197879515
$ws.Range("N18").Value = 1.905094390315568
$ws.Range("B19").Value = 0.4978222760327355
$ws.Range("C19").Value = 0.04496313635235083
$ws.Range("E19").Value = 0.3967552680870696
$ws.Range("F19").Value = 2.549082282000683
$ws.Range("G19").Value = 0.9755393799919148
$ws.Range("H19").Value = 0.9886128695691951
$ws.Range("J19").Value = 0.06986783503639993
$ws.Range("K19").Value = 0.4573019282072721
$ws.Range("M19").Value = 0.4587492710681929
$ws.Range("N19").Value = 1.906662323174778
$ws.Range("B20").Value = 0.5120170362501995
$ws.Range("C20").Value = 0.046652636952075
$ws.Range("E20").Value = 0.4000318853334832
$ws.Range("F20").Value = 2.559014957929037
$ws.Range("G20").Value = 0.9767707294244303
$ws.Range("H20").Value = 0.9875743241342576
$ws.Range("J20").Value = 0.06964436922669126
$ws.Range("K20").Value = 0.4719947354377325
$ws.Range("M20").Value = 0.4666825444994913
$ws.Range("N20").Value = 1.899650229449175
$ws.Range("B21").Value = 0.5598909877232927
$ws.Range("C21").Value = 0.05231926123285291
$ws.Range("E21").Value = 0.4112408867311004
$ws.Range("F21").Value = 2.593940515929518
$ws.Range("G21").Value = 0.9816110602271948
$ws.Range("H21").Value = 0.9845941871860617
$ws.Range("J21").Value = 0.0689180127933966
$ws.Range("K21").Value = 0.5214838391077876
$ws.Range("M21").Value = 0.4935696485959369
$ws.Range("N21").Value = 1.876857350718243
$ws.Range("B22").Value = 0.5913134424070279
$ws.Range("C22").Value = 0.05601534546484288
$ws.Range("E22").Value = 0.4187153836645194
$ws.Range("F22").Value = 2.617924616029086
$ws.Range("G22").Value = 0.9852998730616775
$ws.Range("H22").Value = 0.9830261693458908
$ws.Range("J22").Value = 0.06846152960196639
$ws.Range("K22").Value = 0.5539186832099858
$ws.Range("M22").Value = 0.5113145199879483
$ws.Range("N22").Value = 1.862533207687907
$ws.Range("B23").Value = 0.5745304984561415
$ws.Range("C23").Value = 0.05404335305323116
$ws.Range("E23").Value = 0.4147125115637991
$ws.Range("F23").Value = 2.605017829492809
$ws.Range("G23").Value = 0.9832829287117875
$ws.Range("H23").Value = 0.9838283009168407
$ws.Range("J23").Value = 0.06870349471851256
$ws.Range("K23").Value = 0.5365993494510519
$ws.Range("M23").Value = 0.5018279948675683
$ws.Range("N23").Value = 1.870125851949069
$ws.Range("B24").Value = 0.5112368419228801
$ws.Range("C24").Value = 0.04655989469465283
$ws.Range("E24").Value = 0.3998511959133069
$ws.Range("F24").Value = 2.558463661564943
$ws.Range("G24").Value = 0.9767004675605051
$ws.Range("H24").Value = 0.9876294406647759
$ws.Range("J24").Value = 0.06965654889823503
$ws.Range("K24").Value = 0.4711874102341369
$ws.Range("M24").Value = 0.466246011201541
$ws.Range("N24").Value = 1.900032417471561
$ws.Range("B25").Value = 0.443547576179725
$ws.Range("C25").Value = 0.03845623725362657
$ws.Range("E25").Value = 0.3844615036858841
$ws.Range("F25").Value = 2.513215630703058
$ws.Range("G25").Value = 0.9718473874011408
$ws.Range("H25").Value = 0.9933587932158474
$ws.Range("J25").Value = 0.07076283122657223
$ws.Range("K25").Value = 0.4010265350616748
$ws.Range("M25").Value = 0.4286098999140151
$ws.Range("N25").Value = 1.934743502863377
